$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.259.19"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.865.45"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06552"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07897"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").Value = "1.869.93"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.193"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6824"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "278.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "30.258.12"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.32%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.369"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").Value = "2.111.55"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.201"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.259"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.952"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.383"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09848"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.23%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.482"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.073"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04756"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("E35").Value = "  +4.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7052"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.705"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01879"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.628"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.304"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.959"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8537"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4184"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.228"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "950.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.244"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05646"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.12%  "
